# Cheetah / BaseInit tracking-number refresh.
# Updates the "PackageTrackNum" (col C) / "ShipmentTrackNum" (col D) test
# data on Sheet1 (rows 2-22) with a fresh batch of FedEx tracking numbers.
# Values must stay text (shared-string) cells, matching the original file,
# so we force text entry with a leading apostrophe and then restore the
# "Normal" style so no numFmt/quotePrefix styling sticks to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TrackNum($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

Set-TrackNum "C2"  "320017963792"
Set-TrackNum "C3"  "320017963807"
Set-TrackNum "C4"  "320017963830"
Set-TrackNum "C5"  "320017963851"
Set-TrackNum "D5"  "320017963851"
Set-TrackNum "C6"  "320017963895"
Set-TrackNum "D6"  "320017963895"
Set-TrackNum "C7"  "320017963910"
Set-TrackNum "D7"  "320017963910"
Set-TrackNum "C8"  "320017963943"
Set-TrackNum "C9"  "320017963965"
Set-TrackNum "C10" "320017963998"
Set-TrackNum "C11" "320017964012"
Set-TrackNum "C12" "320017964056"
Set-TrackNum "C13" "320017964078"
Set-TrackNum "D13" "320017964078"
Set-TrackNum "C14" "320017964104"
Set-TrackNum "D14" "320017964104"
Set-TrackNum "C15" "320017964126"
Set-TrackNum "D15" "320017964126"
Set-TrackNum "C16" "320017964159"
Set-TrackNum "D16" "320017964159"
Set-TrackNum "C17" "320017964170"
Set-TrackNum "D17" "320017964170"
Set-TrackNum "C18" "320017964218"
Set-TrackNum "C19" "320017964230"
Set-TrackNum "C20" "320017964284"
Set-TrackNum "C21" "320017964300"
Set-TrackNum "C22" "320017964332"
